$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the worksheet from "login" to "login_failed"
$ws.Name = "login_failed"

# Remove the obsolete "user10" test-data row (old row 2); this shifts the
# remaining rows up and prunes the now-unused shared strings automatically.
$ws.Rows.Item(2).Delete()

# Update the selected cell to reflect the new layout
$ws.Range("C3").Select() | Out-Null

# Update the header/footer font from "Altalanos" to "Regular"
$ps = $ws.PageSetup
$ps.CenterHeader = '&"Times New Roman,Regular"&12&A'
$ps.CenterFooter = '&"Times New Roman,Regular"&12Page &P'
